$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 117, shifting existing rows 117:170 down to 118:171
$ws.Rows.Item(117).Insert()

# Fill the new row 117 with data (copying the pattern of the rest of the table,
# with new values for this entry)
$ws.Cells.Item(117, 1).Value = 5
$ws.Cells.Item(117, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(117, 3).Value = "Maule"
$ws.Cells.Item(117, 4).Value = 44510
$ws.Cells.Item(117, 5).Value = 7
$ws.Cells.Item(117, 6).Value = "Fruta"
$ws.Cells.Item(117, 7).Value = 100108
$ws.Cells.Item(117, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(117, 9).Value = 100108005
$ws.Cells.Item(117, 10).Value = "Piña"
$ws.Cells.Item(117, 11).Value = "Caramelo"
$ws.Cells.Item(117, 12).Value = "Segunda"
$ws.Cells.Item(117, 13).Value = 200
$ws.Cells.Item(117, 14).Value = 18000
$ws.Cells.Item(117, 15).Value = 18000
$ws.Cells.Item(117, 16).Value = 18000
$ws.Cells.Item(117, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(117, 18).Value = "Ecuador"
$ws.Cells.Item(117, 19).Value = 1286
$ws.Cells.Item(117, 20).Value = 14

# Apply date style to new D117 cell (copy from D118, which retains the date style
# after the insert shifted the old D117 down)
$ws.Cells.Item(117, 4).NumberFormat = $ws.Cells.Item(118, 4).NumberFormat
